$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing old rows 21-30 down to 22-31
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new record's data.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T are identical across all rows in this block.
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44438
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108007
$ws.Range("J21").Value = "Coco"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 25
$ws.Range("N21").Value = 21000
$ws.Range("O21").Value = 21000
$ws.Range("P21").Value = 21000
$ws.Range("Q21").Value = "$/malla 20 unidades"
$ws.Range("R21").Value = "Perú"
$ws.Range("S21").Value = 1050
$ws.Range("T21").Value = 20
